# Updated capital structure database
# - Corrects company-name assignment for rows 3 & 4 (Investment Friends SE / Investment
#   Friends Capital SE were swapped) and refreshes the dependent metric columns for the
#   three Estonia "Investments & Asset Management" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix company names: row 3 <-> row 4 ---
$ws.Range("B3").Value = "Investment Friends SE (WSE:IFR)"
$ws.Range("B4").Value = "Investment Friends Capital SE (WSE:IFC)"

# --- Remove stale metrics that no longer apply ---
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("W2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("AA2").ClearContents()
$ws.Range("AC2").ClearContents()
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AA3").ClearContents()
$ws.Range("AC3").ClearContents()
$ws.Range("AI3").ClearContents()
$ws.Range("AK3").ClearContents()
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
$ws.Range("AQ3").ClearContents()

$ws.Range("D4").ClearContents()

# --- Refresh recalculated metrics ---
# Row 2
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1.003448275862069
$ws.Range("J2").Value = 1.003448275862069
$ws.Range("K2").Value = -0.08800000000000002
$ws.Range("L2").Value = -0.303448275862069
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("S2").Value = 0
$ws.Range("U2").Value = 0.014
$ws.Range("V2").Value = 0.003131991051454138
$ws.Range("X2").Value = 0.0389916942384141
$ws.Range("Z2").Value = 0.05331862474719618
$ws.Range("AB2").Value = 0.0389916942384141
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -0.014
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.003141831238779174
$ws.Range("AK2").Value = -0.002604166666666667
$ws.Range("AL2").Value = 0.171
$ws.Range("AM2").Value = 0.03000000000000003
$ws.Range("AO2").Value = 1.701754385964912
$ws.Range("AQ2").Value = 9.699999999999992

# Row 3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1.090322580645161
$ws.Range("J3").Value = 1.090322580645161
$ws.Range("K3").Value = 0.143
$ws.Range("L3").Value = 0.9225806451612902
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("S3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("X3").Value = 0.0389916942384141
$ws.Range("AB3").Value = 0.0389916942384141
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AJ3").Value = 0
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0

# Row 4
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.9037037037037037
$ws.Range("J4").Value = 0.9037037037037037
$ws.Range("K4").Value = -0.231
$ws.Range("L4").Value = -1.711111111111111
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 0.014
$ws.Range("V4").Value = 0.00417910447761194
$ws.Range("W4").Value = -0.04223034734917733
$ws.Range("X4").Value = 0.0389916942384141
$ws.Range("Y4").Value = -0.08122204158759144
$ws.Range("Z4").Value = 0.02482073910645339
$ws.Range("AA4").Value = 0.02243059385916529
$ws.Range("AB4").Value = 0.0389916942384141
$ws.Range("AC4").Value = -0.01656110037924882
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = -0.014
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.004196642685851319
$ws.Range("AK4").Value = -0.002604166666666667
$ws.Range("AL4").Value = 0.171
$ws.Range("AM4").Value = 0.03000000000000003
$ws.Range("AO4").Value = 0.7134502923976608
$ws.Range("AQ4").Value = 4.066666666666663
